# Update the "-1" (column F) and "0" (column G) proportion columns in
# table_S11_family_stats_dn1to0 to reflect the re-run analysis (medians
# instead of median, intra-/inter-group variation) described in the
# commit message. Only rows 2-52 (family entries) are affected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.1565201397829686
$ws.Range("G2").Value = 0.001103549751701306
$ws.Range("F3").Value = 0.4638587456317823
$ws.Range("G3").Value = 0.0007356998344675372
$ws.Range("F4").Value = 0.001103549751701306
$ws.Range("G4").Value = 0
$ws.Range("F5").Value = 0.0003678499172337686
$ws.Range("G5").Value = 0.6864079455582123
$ws.Range("F6").Value = 0.01839249586168843
$ws.Range("G6").Value = 0
$ws.Range("F7").Value = 0.0009196247930844216
$ws.Range("F8").Value = 0.01710502115137024
$ws.Range("G8").Value = 0
$ws.Range("F9").Value = 0.002207099503402612
$ws.Range("G9").Value = 0
$ws.Range("F10").Value = 0.02335846974434431
$ws.Range("G10").Value = 0.001471399668935074
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0.0152657715652014
$ws.Range("F12").Value = 0.04598123965422108
$ws.Range("G12").Value = 0.1329777450800074
$ws.Range("F13").Value = 0.008092698179142909
$ws.Range("G13").Value = 0
$ws.Range("F14").Value = 0.004598123965422108
$ws.Range("G14").Value = 0
$ws.Range("F15").Value = 0.003494574213720802
$ws.Range("G15").Value = 0.000551774875850653
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("F18").Value = 0.000551774875850653
$ws.Range("G18").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("F22").Value = 0.000551774875850653
$ws.Range("G22").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 0.05554533750229906
$ws.Range("G30").Value = 0.001103549751701306
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 0
